$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value() = 131076706
$ws.Range("B6").Value() = 57884
$ws.Range("D6").Value() = "NT"
$ws.Range("E6").Value() = 100109
$ws.Range("F6").Value() = "Tretåig hackspett"
$ws.Range("G6").Value() = "Picoides tridactylus"
$ws.Range("H6").Value() = "(Linnaeus, 1758)"
$ws.Range("P6").Value() = "Stor-Grundsviken, Stor-Grundsviken, Jmt"
$ws.Range("Q6").Value() = 410740
$ws.Range("R6").Value() = 7037755
$ws.Range("S6").Value() = 10
$ws.Range("T6").Value() = "Jämtland"
$ws.Range("U6").Value() = "Åre"
$ws.Range("V6").Value() = "Jämtland"
$ws.Range("W6").Value() = "Kall"

$ws.Range("Y6").Value() = "'2026-02-08"
$ws.Range("Y6").Style = "Normal"
$ws.Range("Z6").Value() = "15:26"
$ws.Range("AA6").Value() = "'2026-02-08"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AB6").Value() = "15:26"

$ws.Range("AD6").Value() = $false
$ws.Range("AE6").Value() = $false
$ws.Range("AG6").Value() = $false

$ws.Range("AW6").Value() = "Fabian Pettersson"
$ws.Range("AX6").Value() = "Fabian Pettersson"
